$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 121.7032956666666
$ws.Range("H2").Value = 365.109887
$ws.Range("I2").Value = 0.8275728186333362
$ws.Range("J2").Value = 0.8275728186333362
$ws.Range("M2").Value = 19.51877066666666
$ws.Range("N2").Value = 58.556312
$ws.Range("O2").Value = 0.9213325069349741
$ws.Range("P2").Value = 0.9213325069349743
$ws.Range("Q2").Value = 2375.498717495193
$ws.Range("R2").Value = 21379.48845745674
$ws.Range("S2").Value = 0.7624697396626944
$ws.Range("T2").Value = 0.7624697396626944
$ws.Range("G3").Value = 121.7032956666666
$ws.Range("H3").Value = 365.109887
$ws.Range("I3").Value = 0.8275728186333362
$ws.Range("J3").Value = 0.8275728186333362
$ws.Range("O3").Value = 0.01004540680524951
$ws.Range("P3").Value = 0.01004540680524951
$ws.Range("Q3").Value = 25.90036800283211
$ws.Range("R3").Value = 233.103312025489
$ws.Range("S3").Value = 0.008313305624138833
$ws.Range("T3").Value = 0.008313305624138833
$ws.Range("G4").Value = 121.7032956666666
$ws.Range("H4").Value = 365.109887
$ws.Range("I4").Value = 0.8275728186333362
$ws.Range("J4").Value = 0.8275728186333362
$ws.Range("N4").Value = 4.361353
$ws.Range("O4").Value = 0.06862208625977624
$ws.Range("P4").Value = 0.06862208625977625
$ws.Range("Q4").Value = 176.9303445552345
$ws.Range("R4").Value = 1592.373100997111
$ws.Range("S4").Value = 0.05678977334650295
$ws.Range("T4").Value = 0.05678977334650297
$ws.Range("I5").Value = 0.1584547531793621
$ws.Range("J5").Value = 0.1584547531793621
$ws.Range("M5").Value = 19.51877066666666
$ws.Range("N5").Value = 58.556312
$ws.Range("O5").Value = 0.9213325069349741
$ws.Range("P5").Value = 0.9213325069349743
$ws.Range("Q5").Value = 454.8349758275031
$ws.Range("R5").Value = 4093.514782447528
$ws.Range("S5").Value = 0.1459895149825043
$ws.Range("T5").Value = 0.1459895149825043
$ws.Range("I6").Value = 0.1584547531793621
$ws.Range("J6").Value = 0.1584547531793621
$ws.Range("O6").Value = 0.01004540680524951
$ws.Range("P6").Value = 0.01004540680524951
$ws.Range("S6").Value = 0.001591742455912095
$ws.Range("T6").Value = 0.001591742455912095
$ws.Range("I7").Value = 0.1584547531793621
$ws.Range("J7").Value = 0.1584547531793621
$ws.Range("N7").Value = 4.361353
$ws.Range("O7").Value = 0.06862208625977624
$ws.Range("P7").Value = 0.06862208625977625
$ws.Range("Q7").Value = 33.87672171584523
$ws.Range("R7").Value = 304.890495442607
$ws.Range("S7").Value = 0.01087349574094574
$ws.Range("T7").Value = 0.01087349574094574
$ws.Range("I8").Value = 0.01397242818730167
$ws.Range("J8").Value = 0.01397242818730167
$ws.Range("M8").Value = 19.51877066666666
$ws.Range("N8").Value = 58.556312
$ws.Range("O8").Value = 0.9213325069349741
$ws.Range("P8").Value = 0.9213325069349743
$ws.Range("Q8").Value = 40.10702682821511
$ws.Range("R8").Value = 360.963241453936
$ws.Range("S8").Value = 0.01287325228977555
$ws.Range("T8").Value = 0.01287325228977555
$ws.Range("I9").Value = 0.01397242818730167
$ws.Range("J9").Value = 0.01397242818730167
$ws.Range("O9").Value = 0.01004540680524951
$ws.Range("P9").Value = 0.01004540680524951
$ws.Range("S9").Value = 0.0001403587251985803
$ws.Range("T9").Value = 0.0001403587251985803
$ws.Range("I10").Value = 0.01397242818730167
$ws.Range("J10").Value = 0.01397242818730167
$ws.Range("N10").Value = 4.361353
$ws.Range("O10").Value = 0.06862208625977624
$ws.Range("P10").Value = 0.06862208625977625
$ws.Range("S10").Value = 0.0009588171723275445
$ws.Range("T10").Value = 0.0009588171723275445
